# Update the "Förändrad" (Changed) date column (C) for data rows 2-18
# from serial date 45171 (2023-09-02) to 45172 (2023-09-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }
}
